$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text format to A:C for rows 2-6, then write the new bank code values
# into column B as text (the upload format treats bank codes as strings).
$ws.Range("A2:C6").NumberFormat = "@"
$ws.Range("E2:E6").NumberFormat = "0.00"

$ws.Range("B2").Value = "199999"
$ws.Range("B3").Value = "288888"
$ws.Range("B4").Value = "388888"
$ws.Range("B5").Value = "488888"
$ws.Range("B6").Value = "588888"

[void]$ws.Range("C3").Select()

$ws.PageSetup.Orientation = 1
